# Insert a new data row after the existing row 21 (i.e. as the new row 22),
# shifting all subsequent rows (old 22..145) down by one to become 23..146.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with the new weekly record.
# Non date/volume columns mirror the row that used to occupy row 22
# (now shifted to row 23), since those attributes are constant for this
# market/category subset.
$ws.Range("A22").Value = 4
$ws.Range("B22").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C22").Value = "Los Lagos"
$ws.Range("D22").Value = 44462
$ws.Range("E22").Value = 10
$ws.Range("F22").Value = 100112043
$ws.Range("G22").Value = "Pepino ensalada"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 19000
$ws.Range("L22").Value = 19000
$ws.Range("M22").Value = 19000
$ws.Range("N22").Value = "`$/caja 60 unidades"
$ws.Range("O22").Value = "Región de Arica y Parinacota"
$ws.Range("P22").Value = 317
$ws.Range("Q22").Value = 60
$ws.Range("R22").Value = "Hortaliza"
